# Atualização automática de preços de eletricidade
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Day
$ws.Range("A2").Value = 45893

# Hourly prices 0h-1h .. 23h-24h (columns B..Y)
$ws.Range("B2").Value = 106
$ws.Range("C2").Value = 101.12
$ws.Range("D2").Value = 97.52
$ws.Range("E2").Value = 94.23999999999999
$ws.Range("F2").Value = 93.76000000000001
$ws.Range("G2").Value = 94
$ws.Range("H2").Value = 94.23999999999999
$ws.Range("I2").Value = 95.95999999999999
$ws.Range("J2").Value = 75
$ws.Range("K2").Value = 25
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = -0.01
$ws.Range("N2").Value = -0.01
$ws.Range("O2").Value = 0
$ws.Range("P2").Value = -0.01
$ws.Range("Q2").Value = -0.1
$ws.Range("R2").Value = 0
$ws.Range("S2").Value = 12.93
$ws.Range("T2").Value = 61.15
$ws.Range("U2").Value = 98.56999999999999
$ws.Range("V2").Value = 114.99
$ws.Range("W2").Value = 114.32
$ws.Range("X2").Value = 111.92
$ws.Range("Y2").Value = 105.77

# Price_Daily_Avg
$ws.Range("Z2").Value = 62.35

# Slot_4h_max (unchanged) / Slot_4h_price
$ws.Range("AB2").Value = 111.75

# Slot_2h_frist (unchanged) / Slot_2h_frist_price
$ws.Range("AD2").Value = 114.66

# Slot_2h_second (unchanged) / Slot_2h_second_price
$ws.Range("AF2").Value = 108.84

# Slot_min_price
$ws.Range("AG2").Value = "9h-18h"
